# daily auto push: 2026-01-20 22:35 UTC
# Insert a new data row at row 681 (pushing existing rows 681:722 down to
# 682:723) and populate it with the new reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(681).Insert()

$ws.Cells.Item(681, 1).NumberFormat = "@"
$ws.Cells.Item(681, 1).Value = "2026/01/21"
$ws.Cells.Item(681, 2).Value = "水"
$ws.Cells.Item(681, 3).Value = 6
$ws.Cells.Item(681, 4).Value = 201
